$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.350.23'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '2.240.03'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''245.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '''74.41'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.03%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").Value = '''43.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.35%  '
$ws.Range("D11").Value = '''0.0955'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '''0.855'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '2.241.60'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '42.264.99'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("D20").Value = '''72.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = '''10.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +42.55%  '
$ws.Range("D22").Value = '''231.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '''2.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.66%  '
$ws.Range("D24").Value = '''11.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.63%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = '''2.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.44%  '
$ws.Range("D29").Value = '''166.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").Value = '''20.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.63%  '
$ws.Range("D31").Value = '''5.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +20.01%  '
$ws.Range("E32").Value = '  -1.49%  '
$ws.Range("E34").Value = '  -9.06%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '''4.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("E38").Value = '  -5.26%  '
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").Value = '''5.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.00%  '
$ws.Range("D41").Value = '''63.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.09%  '
$ws.Range("D42").Value = '''0.203'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").Value = '''8.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.10%  '
$ws.Range("D44").Value = '''105.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.92%  '
$ws.Range("D45").Value = '''0.103'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.57%  '
$ws.Range("D46").Value = '''0.995'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.40%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''1.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("D51").Value = '''4.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.75%  '
